$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; this shifts the existing rows 19-77 down to 20-78
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new record
$ws.Range("A19").Value = 5
$ws.Range("B19").Value = "Macroferia Regional de Talca"
$ws.Range("C19").Value = "Maule"
$ws.Range("D19").Value = 44802
$ws.Range("E19").Value = 7
$ws.Range("F19").Value = 100112026
$ws.Range("G19").Value = "Haba"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 12000
$ws.Range("M19").Value = 12000
$ws.Range("N19").Value = '$/saco 25 kilos'
$ws.Range("O19").Value = 'Provincia del Elquí'
$ws.Range("P19").Value = 480
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
